# Updates the "cryptos" list with newly scraped Price / Volume(1h) figures.
# Each target cell is stored as an inline string in the workbook, so we force
# the cell to text ("@") before writing the value and then restore the
# default ("Normal") cell style so no stray number formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "27.095.05"
    "E2" = "  -0.09%  "
    "D3" = "1.623.74"
    "E3" = "  -0.92%  "
    "E4" = "  -0.13%  "
    "D5" = "214.36"
    "E5" = "  -1.05%  "
    "E6" = "  -1.00%  "
    "E7" = "  -0.15%  "
    "E8" = "  +0.70%  "
    "E9" = "  -1.48%  "
    "D10" = "20.02"
    "E10" = "  +0.25%  "
    "D11" = "0.0846"
    "D12" = "1.851.33"
    "D13" = "1.620.88"
    "E13" = "  -1.11%  "
    "D14" = "4.14"
    "E14" = "  +0.10%  "
    "E15" = "  +0.04%  "
    "D16" = "27.057.65"
    "E16" = "  -0.29%  "
    "D17" = "64.53"
    "E17" = "  -3.28%  "
    "E18" = "  +0.13%  "
    "D19" = "213.49"
    "E19" = "  -1.49%  "
    "E20" = "  -0.03%  "
    "D21" = "6.86"
    "E21" = "  -1.58%  "
    "D22" = "4.36"
    "E22" = "  -1.25%  "
    "E23" = "  -6.96%  "
    "D24" = "9.06"
    "E24" = "  -0.82%  "
    "D25" = "148.26"
    "E25" = "  +0.93%  "
    "E26" = "  -0.05%  "
    "D27" = "7.37"
    "E27" = "  -0.62%  "
    "E28" = "  -2.88%  "
    "D29" = "15.54"
    "E29" = "  -0.88%  "
    "E30" = "  +0.35%  "
    "E31" = "  -0.82%  "
    "E32" = "  -0.94%  "
    "D33" = "0.734"
    "E33" = "  +35.39%  "
    "E34" = "  -0.46%  "
    "D35" = "1.360.92"
    "E35" = "  +4.13%  "
    "E36" = "  -0.18%  "
    "D37" = "2.45"
    "E37" = "  -0.64%  "
    "D39" = "0.845"
    "E39" = "  -1.45%  "
    "E40" = "  -0.15%  "
    "E41" = "  -1.30%  "
    "E42" = "  +0.50%  "
    "D43" = "64.35"
    "E43" = "  +4.10%  "
    "E44" = "  +0.77%  "
    "D45" = "1.762.40"
    "E45" = "  -0.92%  "
    "E46" = "  +3.52%  "
    "D47" = "89.95"
    "E47" = "  -1.60%  "
    "D48" = "0.871"
    "E48" = "  +30.37%  "
    "E49" = "  -1.81%  "
    "E50" = "  +5.00%  "
    "E51" = "  +0.04%  "
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}
